$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'67.688.27"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "  -1.03%  "

$c = $ws.Range("D3")
$c.Value = "'3.788.93"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "  +0.04%  "

$c = $ws.Range("E4")
$c.Value = "  +0.05%  "

$c = $ws.Range("D5")
$c.Value = "'595.84"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "  +0.05%  "

$c = $ws.Range("D6")
$c.Value = "'166.87"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "  -0.57%  "

$c = $ws.Range("D7")
$c.Value = "'3.785.66"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "  -0.02%  "

$c = $ws.Range("E8")
$c.Value = "  +0.01%  "

$c = $ws.Range("D9")
$c.Value = "'0.520"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "  -0.19%  "

$c = $ws.Range("E10")
$c.Value = "  -0.67%  "

$c = $ws.Range("D11")
$c.Value = "'6.34"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "  -1.17%  "

$c = $ws.Range("E12")
$c.Value = "  -0.42%  "

$c = $ws.Range("D13")
$c.Value = "'0.0000253"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "  -2.58%  "

$c = $ws.Range("D14")
$c.Value = "'36.06"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "  -0.29%  "

$c = $ws.Range("D15")
$c.Value = "'4.420.59"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "  +0.11%  "

$c = $ws.Range("D16")
$c.Value = "'3.773.72"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "  +0.29%  "

$c = $ws.Range("D17")
$c.Value = "'18.53"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "  +3.19%  "

$c = $ws.Range("D18")
$c.Value = "'67.638.91"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "  -1.01%  "

$c = $ws.Range("D19")
$c.Value = "'7.06"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "  +1.19%  "

$c = $ws.Range("E20")
$c.Value = "  +0.05%  "

$c = $ws.Range("E21")
$c.Value = "  -7.73%  "

$c = $ws.Range("D22")
$c.Value = "'459.34"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "  -1.36%  "

$c = $ws.Range("D23")
$c.Value = "'0.699"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "  -0.13%  "

$c = $ws.Range("E24")
$c.Value = "  +3.03%  "

$c = $ws.Range("E25")
$c.Value = "  -0.67%  "

$c = $ws.Range("D26")
$c.Value = "'12.06"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "  +1.30%  "

$c = $ws.Range("E27")
$c.Value = "  -3.14%  "

$c = $ws.Range("E28")
$c.Value = "  -0.06%  "

$c = $ws.Range("E29")
$c.Value = "  -1.61%  "

$c = $ws.Range("D30")
$c.Value = "'3.933.64"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "  +0.07%  "

$c = $ws.Range("E31")
$c.Value = "  +0.04%  "

$c = $ws.Range("E32")
$c.Value = "  +4.12%  "

$c = $ws.Range("D33")
$c.Value = "'7.23"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "  -1.24%  "

$c = $ws.Range("D34")
$c.Value = "'29.59"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "  -1.48%  "

$c = $ws.Range("B35")
$c.Value = "Aptos"
$c = $ws.Range("C35")
$c.Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D35")
$c.Value = "'9.07"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "  -1.05%  "

$c = $ws.Range("B36")
$c.Value = "Binance-PegBSC-USD"
$c = $ws.Range("C36")
$c.Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D36")
$c.Value = "'0.993"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "  -0.61%  "

$c = $ws.Range("D37")
$c.Value = "'0.100"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "  -0.69%  "

$c = $ws.Range("E38")
$c.Value = "  -2.65%  "

$c = $ws.Range("D39")
$c.Value = "'0.138"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "  -0.29%  "

$c = $ws.Range("D40")
$c.Value = "'0.992"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "  -0.88%  "

$c = $ws.Range("D41")
$c.Value = "'5.77"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "  -0.46%  "

$c = $ws.Range("E42")
$c.Value = "  +0.04%  "

$c = $ws.Range("E43")
$c.Value = "  +0.01%  "

$c = $ws.Range("D44")
$c.Value = "'48.01"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "  +2.01%  "

$c = $ws.Range("D45")
$c.Value = "'43.99"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "  +0.02%  "

$c = $ws.Range("E46")
$c.Value = "  -1.43%  "

$c = $ws.Range("D47")
$c.Value = "'150.79"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "  +3.05%  "

$c = $ws.Range("D48")
$c.Value = "'8.28"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "  -1.80%  "

$c = $ws.Range("D49")
$c.Value = "'26.70"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "  +5.66%  "

$c = $ws.Range("D50")
$c.Value = "'388.69"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "  -1.71%  "

$c = $ws.Range("E51")
$c.Value = "  -5.08%  "
